# Update the "Förändrad" (Changed) date column (C) from 45175 to 45177
# for all data rows (rows 2 through 307) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C307").Value = 45177
